$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1. Resize / reposition existing shapes so the big "TaskListPanel" region is
#    split in two, making room for a new "CalendarPanel" to its right.
# ---------------------------------------------------------------------------

# Rectangle 8 (big panel background, first row of panels)
$sh = $s.Shapes.Item(6)
$sh.Left = 245.3792955
$sh.Width = 214.75858350000001

# Rectangle 9 (header bar inside first panel)
$sh = $s.Shapes.Item(7)
$sh.Left = 251.1725235
$sh.Width = 202.758507

# TextBox 10 ("TaskCard" label, first panel)
$sh = $s.Shapes.Item(8)
$sh.Left = 289.7998505

# TextBox 13 ("TaskListPanel" label)
$sh = $s.Shapes.Item(9)
$sh.Left = 289.7998505

# Rectangle 17 (header bar inside second panel)
$sh = $s.Shapes.Item(14)
$sh.Left = 251.1725235
$sh.Width = 202.758507

# TextBox 18 ("TaskCard" label, second panel)
$sh = $s.Shapes.Item(15)
$sh.Left = 289.7998505

# TextBox 11 (". . ." vertical label)
$sh = $s.Shapes.Item(17)
$sh.Left = 341.793076

# ---------------------------------------------------------------------------
# 2. Add the new "CalendarPanel" rectangle + label, duplicated from existing
#    shapes so style/fill/effects match exactly, then repositioned/retexted.
# ---------------------------------------------------------------------------

# Consume the lowest free shape-id slots first so the new shapes land on the
# next ids after the highest existing one (21 / 22), matching the target deck.
$gap1 = $s.Shapes.AddShape(1, 0, 0, 10, 10)
$gap2 = $s.Shapes.AddShape(1, 0, 0, 10, 10)
$gap1.Delete()
$gap2.Delete()

$panelSrc = $s.Shapes.Item(6)
$panelDup = $panelSrc.Duplicate()
$newPanel = $panelDup.Item(1)
$newPanel.Name = "Rectangle 20"
$newPanel.Left = 468.413773
$newPanel.Top = 203.80953250000002
$newPanel.Width = 214.75858350000001
$newPanel.Height = 218.89653800000002

$labelSrc = $s.Shapes.Item(8)
$labelDup = $labelSrc.Duplicate()
$newLabel = $labelDup.Item(1)
$newLabel.Name = "TextBox 21"
$newLabel.Left = 509.37928800000003
$newLabel.Top = 291.6705475
$newLabel.TextFrame.TextRange.Text = "CalendarPanel"

# ---------------------------------------------------------------------------
# 3. Refresh the cached "datetimeFigureOut" date text on the master & every
#    layout (best-effort text update; value rolls from 10/10/2016 -> 10/25/2016).
# ---------------------------------------------------------------------------

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $msh = $master.Shapes.Item($i)
    if ($msh.Name -like "Date Placeholder*") {
        $msh.TextFrame.TextRange.Text = "10/25/2016"
    }
}

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $cl = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $lsh = $cl.Shapes.Item($i)
        if ($lsh.Name -like "Date Placeholder*") {
            $lsh.TextFrame.TextRange.Text = "10/25/2016"
        }
    }
}
